$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D header and first two tag values (reused later)
$ws.Range("D1").Value = "tag_id"
$ws.Range("D2").Value = "tag_1"
$ws.Range("D3").Value = "tag_3"

# Append new rows 4-8
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "alice 2"
$ws.Range("C4").Value = 151
$ws.Range("D4").Value = "tag_1"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "bob 2"
$ws.Range("C5").Value = 210
$ws.Range("D5").Value = "tag_1"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "alice 3"
$ws.Range("C6").Value = 269
$ws.Range("D6").Value = "tag_1"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "bob 3"
$ws.Range("C7").Value = 328
$ws.Range("D7").Value = "tag_3"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "alice 3"
$ws.Range("C8").Value = 387
$ws.Range("D8").Value = "tag_1"

# Rename existing alice/bob to alice 1 / bob 1
$ws.Range("B2").Value = "alice 1"
$ws.Range("B3").Value = "bob 1"

# Match the final selection state (entire row 9 selected, as if ready for new entry)
[void]$ws.Rows(9).Select()
